$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

# --- Footer (default/primary) -> word/footer2.xml: Pearson logo docPr id="2" ---
$ftr1 = $sec.Footers.Item(1)
$shpF1 = $ftr1.Range.InlineShapes.Item(1)
$shpF1.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# --- Footer (first page) -> word/footer1.xml: Pearson logo docPr id="3" ---
$ftr2 = $sec.Footers.Item(2)
$shpF2 = $ftr2.Range.InlineShapes.Item(1)
$shpF2.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image2.png"

# --- Header (first page) -> word/header1.xml: BTec logo docPr id="1" ---
$hdr2 = $sec.Headers.Item(2)
$shpH2 = $hdr2.Range.InlineShapes.Item(1)
$shpH2.Range.Select()
$word.Selection.InlineShapes.Item(1).Name = "image1.jpg"
